# Apply edits to the "articles" worksheet of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("articles")

# Update the material number codes (column A) for rows 7-10.
# B column (material names) stay unchanged.
$ws.Range("A7").Value = "000011"
$ws.Range("A8").Value = "000012"
$ws.Range("A9").Value = "000012"
$ws.Range("A10").Value = "000203"

# Update the active selection on the sheet to reflect where the user left off.
$ws.Range("I10").Select()
